$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.78"
$ws.Range("E2").Value = "'0.19%"
$ws.Range("G2").Value = "'16"
$ws.Range("D3").Value = "'43.99"
$ws.Range("E3").Value = "'0.28%"
$ws.Range("G3").Value = "'16"
$ws.Range("D4").Value = "'5.482"
$ws.Range("E4").Value = "'-1.48%"
$ws.Range("G4").Value = "'16"
$ws.Range("D5").Value = "'0.08034"
$ws.Range("E5").Value = "'-0.64%"
$ws.Range("G5").Value = "'16"
$ws.Range("D6").Value = "'2.034"
$ws.Range("E6").Value = "'6.39%"
$ws.Range("G6").Value = "'16"
$ws.Range("D7").Value = "'0.9529"
$ws.Range("E7").Value = "'0.25%"
$ws.Range("G7").Value = "'16"
$ws.Range("D8").Value = "'0.1103"
$ws.Range("E8").Value = "'-8.25%"
$ws.Range("G8").Value = "'16"
$ws.Range("D9").Value = "'0.1874"
$ws.Range("E9").Value = "'1.12%"
$ws.Range("G9").Value = "'16"
$ws.Range("D10").Value = "'10.16"
$ws.Range("E10").Value = "'2.19%"
$ws.Range("G10").Value = "'16"
$ws.Range("D11").Value = "'0.09995"
$ws.Range("E11").Value = "'2.91%"
$ws.Range("G11").Value = "'16"
$ws.Range("D12").Value = "'0.04738"
$ws.Range("E12").Value = "'7.70%"
$ws.Range("G12").Value = "'16"
$ws.Range("D13").Value = "'0.1058"
$ws.Range("E13").Value = "'-0.84%"
$ws.Range("G13").Value = "'16"
$ws.Range("D14").Value = "'0.001269"
$ws.Range("E14").Value = "'-0.23%"
$ws.Range("G14").Value = "'16"
$ws.Range("D15").Value = "'0.04077"
$ws.Range("E15").Value = "'-2.62%"
$ws.Range("G15").Value = "'16"
$ws.Range("D16").Value = "'0.005808"
$ws.Range("E16").Value = "'-2.15%"
$ws.Range("G16").Value = "'16"
$ws.Range("D17").Value = "'3.372"
$ws.Range("E17").Value = "'-0.92%"
$ws.Range("G17").Value = "'16"
$ws.Range("D18").Value = "'4.413"
$ws.Range("E18").Value = "'3.25%"
$ws.Range("G18").Value = "'16"
$ws.Range("D19").Value = "'2.613"
$ws.Range("E19").Value = "'2.51%"
$ws.Range("G19").Value = "'16"
$ws.Range("D20").Value = "'0.3414"
$ws.Range("E20").Value = "'-0.44%"
$ws.Range("G20").Value = "'16"
$ws.Range("E21").Value = "'-0.62%"
$ws.Range("G21").Value = "'16"
$ws.Range("D22").Value = "'0.2575"
$ws.Range("E22").Value = "'3.17%"
$ws.Range("G22").Value = "'16"
$ws.Range("D23").Value = "'0.001309"
$ws.Range("E23").Value = "'5.59%"
$ws.Range("G23").Value = "'16"
$ws.Range("D24").Value = "'0.004338"
$ws.Range("E24").Value = "'-0.07%"
$ws.Range("G24").Value = "'16"
$ws.Range("D25").Value = "'0.0001251"
$ws.Range("E25").Value = "'5.42%"
$ws.Range("G25").Value = "'16"
$ws.Range("D26").Value = "'0.0003742"
$ws.Range("E26").Value = "'-5.83%"
$ws.Range("G26").Value = "'16"
$ws.Range("G27").Value = "'16"
$ws.Range("G28").Value = "'16"
$ws.Range("G29").Value = "'16"
$ws.Range("G30").Value = "'16"
$ws.Range("G31").Value = "'16"
$ws.Range("G32").Value = "'16"
$ws.Range("G33").Value = "'16"
$ws.Range("G34").Value = "'16"
$ws.Range("G35").Value = "'16"
$ws.Range("G36").Value = "'16"
$ws.Range("G37").Value = "'16"
$ws.Range("D38").Value = "'0.02577"
$ws.Range("E38").Value = "'-3.95%"
$ws.Range("G38").Value = "'16"
$ws.Range("D39").Value = "'0.05683"
$ws.Range("E39").Value = "'2.93%"
$ws.Range("G39").Value = "'16"
$ws.Range("D40").Value = "'0.007652"
$ws.Range("E40").Value = "'0.78%"
$ws.Range("G40").Value = "'16"
$ws.Range("D41").Value = "'0.1397"
$ws.Range("E41").Value = "'-0.53%"
$ws.Range("G41").Value = "'16"
$ws.Range("D42").Value = "'0.007395"
$ws.Range("E42").Value = "'-8.93%"
$ws.Range("G42").Value = "'16"
$ws.Range("D43").Value = "'0.002009"
$ws.Range("E43").Value = "'0.74%"
$ws.Range("G43").Value = "'16"
$ws.Range("D44").Value = "'0.008504"
$ws.Range("E44").Value = "'-3.15%"
$ws.Range("G44").Value = "'16"
$ws.Range("D45").Value = "'0.00007100"
$ws.Range("E45").Value = "'0.66%"
$ws.Range("G45").Value = "'16"
$ws.Range("E46").Value = "'0.44%"
$ws.Range("G46").Value = "'16"
$ws.Range("D47").Value = "'0.0005810"
$ws.Range("E47").Value = "'-0.02%"
$ws.Range("G47").Value = "'16"
$ws.Range("D48").Value = "'0.003500"
$ws.Range("E48").Value = "'54.79%"
$ws.Range("G48").Value = "'16"
$ws.Range("D49").Value = "'0.003508"
$ws.Range("E49").Value = "'7.68%"
$ws.Range("G49").Value = "'16"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.44%"
$ws.Range("G50").Value = "'16"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.44%"
$ws.Range("G51").Value = "'16"
